# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 46074 (2026-02-21) to 46075 (2026-02-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 128 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 46075
    }
}
